# Refresh the crypto price/volume snapshot (and two swapped row pairs) per the
# upstream scraper run. Numeric-looking Price values are written with a leading
# apostrophe so Excel keeps them as text (preserving trailing zeros / exact
# digits, e.g. '1.00' or '0.0000245') instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.865.45'
$ws.Range('E2').Value = '  +0.14%  '

$ws.Range('D3').Value = '3.788.85'
$ws.Range('E3').Value = '  -0.45%  '

$ws.Range('D4').Value = '''0.997'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '''600.32'
$ws.Range('E5').Value = '  -0.89%  '

$ws.Range('D6').Value = '''163.29'
$ws.Range('E6').Value = '  -2.24%  '

$ws.Range('D7').Value = '3.787.55'
$ws.Range('E7').Value = '  -0.48%  '

$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.29%  '

$ws.Range('D9').Value = '''0.513'
$ws.Range('E9').Value = '  -1.42%  '

$ws.Range('D10').Value = '''0.157'
$ws.Range('E10').Value = '  -2.70%  '

$ws.Range('E11').Value = '  -0.88%  '

$ws.Range('D12').Value = '''6.62'
$ws.Range('E12').Value = '  +5.13%  '

$ws.Range('D13').Value = '''0.0000245'
$ws.Range('E13').Value = '  -3.54%  '

$ws.Range('D14').Value = '''35.16'
$ws.Range('E14').Value = '  -2.24%  '

$ws.Range('D15').Value = '4.422.69'
$ws.Range('E15').Value = '  -0.48%  '

$ws.Range('D16').Value = '3.778.08'
$ws.Range('E16').Value = '  -1.23%  '

$ws.Range('D17').Value = '67.825.79'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').Value = '''18.14'
$ws.Range('E18').Value = '  -1.85%  '

$ws.Range('E19').Value = '  +2.06%  '

$ws.Range('D20').Value = '''7.01'
$ws.Range('E20').Value = '  -1.09%  '

$ws.Range('D21').Value = '''458.01'
$ws.Range('E21').Value = '  -0.91%  '

$ws.Range('D22').Value = '''9.46'
$ws.Range('E22').Value = '  -4.46%  '

$ws.Range('E23').Value = '  -1.28%  '

$ws.Range('D24').Value = '''83.04'
$ws.Range('E24').Value = '  -0.40%  '

$ws.Range('D25').Value = '''0.0000143'
$ws.Range('E25').Value = '  -5.76%  '

$ws.Range('D26').Value = '''11.86'
$ws.Range('E26').Value = '  -1.86%  '

$ws.Range('E27').Value = '  -1.10%  '

$ws.Range('D29').Value = '''9.91'
$ws.Range('E29').Value = '  -0.99%  '

$ws.Range('D30').Value = '3.933.19'
$ws.Range('E30').Value = '  -0.56%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''7.21'
$ws.Range('E31').Value = '  -2.68%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''2.19'
$ws.Range('E32').Value = '  -1.94%  '

$ws.Range('E33').Value = '  -7.49%  '

$ws.Range('D34').Value = '''28.91'
$ws.Range('E34').Value = '  -2.28%  '

$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('E36').Value = '  -1.48%  '

$ws.Range('D37').Value = '''0.0990'
$ws.Range('E37').Value = '  -0.94%  '

$ws.Range('E38').Value = '  +5.16%  '

$ws.Range('D39').Value = '''5.81'
$ws.Range('E39').Value = '  +0.43%  '

$ws.Range('D40').Value = '''0.980'
$ws.Range('E40').Value = '  -2.05%  '

$ws.Range('D41').Value = '''3.17'
$ws.Range('E41').Value = '  -6.52%  '

$ws.Range('E42').Value = '  -0.02%  '

$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('D44').Value = '''44.01'
$ws.Range('E44').Value = '  +2.18%  '

$ws.Range('D45').Value = '''47.14'
$ws.Range('E45').Value = '  -2.14%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '''0.295'
$ws.Range('E46').Value = '  -2.23%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '''151.82'
$ws.Range('E47').Value = '  +2.18%  '

$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = '''8.28'
$ws.Range('E48').Value = '  -0.87%  '

$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '''1.35'
$ws.Range('E49').Value = '  -1.26%  '

$ws.Range('E50').Value = '  -0.58%  '

$ws.Range('D51').Value = '''383.09'
$ws.Range('E51').Value = '  -2.16%  '
